$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 40, shifting existing rows 40-46 down to 41-47
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new weekly record
$ws.Cells.Item(40, 1).Value = 4
$ws.Cells.Item(40, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(40, 3).Value = "Los Lagos"
$ws.Cells.Item(40, 4).Value = 45209
$ws.Cells.Item(40, 5).Value = 10
$ws.Cells.Item(40, 6).Value = 100112013
$ws.Cells.Item(40, 7).Value = "Alcachofa"
$ws.Cells.Item(40, 8).Value = "Española"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 150
$ws.Cells.Item(40, 11).Value = 12000
$ws.Cells.Item(40, 12).Value = 12000
$ws.Cells.Item(40, 13).Value = 12000
$ws.Cells.Item(40, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(40, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(40, 16).Value = 400
$ws.Cells.Item(40, 17).Value = 30
$ws.Cells.Item(40, 18).Value = "Hortaliza"
